$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$excel.ActiveWindow.ScrollColumn = [int]19
$excel.ActiveWindow.ScrollRow = [int]1
$ws.Range("Z21").Select()
Write-Output "done"
